# Insert a new weekly price-report row right after the existing row 36
# (i.e. at row 37), pushing all the subsequent "Espinaca" records down by
# one row (old row 37 -> new row 38, ..., old row 204 -> new row 205).
#
# The new row 37 carries a new observation for 2021-12-31 (serial 44561)
# with Volumen (J) = 2900; every other column for that row keeps the same
# values the series already used (they are the constant/static columns
# for this market+category), so we clone them from the row that lands
# right below it after the insert (new row 38, i.e. the old row 37) and
# then overwrite just D37/J37.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 37..204 down to 38..205, leaving a blank row 37 behind.
$ws.Rows.Item(37).Insert()

# Clone the (now shifted) old row 37 data - currently sitting at row 38 -
# into the freshly inserted blank row 37.
$ws.Rows.Item(38).Copy()
$ws.Rows.Item(37).PasteSpecial()

# Overwrite the new row's date and volume with the new data point.
$ws.Range("D37").Value = 44561
$ws.Range("J37").Value = 2900
